$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.618329644203186
$ws.Range("B1").Value = 2.563636541366577
$ws.Range("C1").Value = 2.819139957427979
$ws.Range("D1").Value = 3.158385753631592
$ws.Range("E1").Value = 3.162692070007324
